$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A71").Value = "Matrix-Multiplication"
$ws.Range("B71").Value = "矩阵连乘问题——算法设计课程"
